# Generate Report for Handoff
#
# The 2acdd304-... file has just been handed off again (new handoff round),
# so it now sorts ahead of 929aaf6c-... everywhere it is listed, and its
# status flips from "Handed back: in sync with en-US" to "Ready for handoff".
# The 929aaf6c-... file stays as it was (still "Handed back: in sync with en-US").
#
# Concretely: on every sheet, the row that used to show the 2acdd304 file
# (row 2) and the row that used to show the 929aaf6c file (row 3) swap
# places, and the (now row 3) 2acdd304 row gets the refreshed status /
# handoff timestamp.

$wb = $excel.ActiveWorkbook

$uuidA = "2acdd304-ba47-4fb6-9be2-0e5e51f2233c"
$uuidB = "929aaf6c-2192-4169-ac76-63d9b8244b10"

$readyForHandoff = "Ready for handoff"
$handedBack = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value2 = ($uuidB + ".md")
$wsOverview.Range("B2").Value2 = $handedBack
$wsOverview.Range("C2").Value2 = $handedBack

$wsOverview.Range("A3").Value2 = ($uuidA + ".md")
$wsOverview.Range("B3").Value2 = $readyForHandoff
$wsOverview.Range("C3").Value2 = $readyForHandoff

$wsOverview.Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/38beca1a811b8e315158bc4dda8d495b6acea141/e2e/" + $uuidB + ".md", "", "", $uuidB + ".md")
$wsOverview.Hyperlinks.Add($wsOverview.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/38beca1a811b8e315158bc4dda8d495b6acea141/e2e/" + $uuidA + ".md", "", "", $uuidA + ".md")
$wsOverview.Hyperlinks.Add($wsOverview.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/38beca1a811b8e315158bc4dda8d495b6acea141/.localization-config", "", "", ".localization-config")

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$zhHandoffCommit = "de042f8f5b7dd8f1e493ce48ccb31c7e9a5be108"
$zhHandoffRepoCommitA = "8fdd61e8305cac1bb1fd5f124854164bc7358e7b"
$zhHandbackCommit = "abfc6699284cec51ff0900e5cab0d451d30ce071"
$zhHash = "7d99b0cfb94d07c6c9898a5664a9f07e6b3c9ea5"
$zhHashB = "fa8dfbcb6f034a590a88c3ecb061e2853a9c33f6"

# Row 2 now holds the 929aaf6c file (status/timestamps unchanged from before)
$wsZhCn.Range("A2").Value2 = ($uuidB + ".md")
$wsZhCn.Range("B2").Value2 = $handedBack
$wsZhCn.Range("C2").Value2 = ($uuidB + "." + $zhHashB + ".zh-cn.xlf")
$wsZhCn.Range("D2").Value2 = "2016-03-08 18:49:04"
$wsZhCn.Range("E2").Value2 = ($uuidB + ".md")
$wsZhCn.Range("F2").Value2 = ($uuidB + "." + $zhHashB + ".zh-cn.xlf")
$wsZhCn.Range("G2").Value2 = "2016-03-08 18:50:14"
$wsZhCn.Range("H2").Value2 = "Include"

# Row 3 now holds the 2acdd304 file, with refreshed status + handoff datetime
$wsZhCn.Range("A3").Value2 = ($uuidA + ".md")
$wsZhCn.Range("B3").Value2 = $readyForHandoff
$wsZhCn.Range("C3").Value2 = ($uuidA + "." + $zhHash + ".zh-cn.xlf")
$wsZhCn.Range("D3").Value2 = "2016-03-08 18:51:09"
$wsZhCn.Range("E3").Value2 = ($uuidA + ".md")
$wsZhCn.Range("F3").Value2 = ($uuidA + "." + $zhHash + ".zh-cn.xlf")
$wsZhCn.Range("G3").Value2 = "2016-03-08 18:50:14"
$wsZhCn.Range("H3").Value2 = "Include"

$wsZhCn.Hyperlinks.Delete()
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/38beca1a811b8e315158bc4dda8d495b6acea141/e2e/" + $uuidB + ".md", "", "", $uuidB + ".md")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/" + $zhHandoffCommit + "/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/" + $uuidB + "." + $zhHashB + ".zh-cn.xlf", "", "", $uuidB + "." + $zhHashB + ".zh-cn.xlf")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("E2"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/" + $zhHandoffRepoCommitA + "/e2e/" + $uuidB + ".md", "", "", $uuidB + ".md")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/" + $zhHandbackCommit + "/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/" + $uuidB + "." + $zhHashB + ".zh-cn.xlf", "", "", $uuidB + "." + $zhHashB + ".zh-cn.xlf")

$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/38beca1a811b8e315158bc4dda8d495b6acea141/e2e/" + $uuidA + ".md", "", "", $uuidA + ".md")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("C3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/" + $zhHandoffCommit + "/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/" + $uuidA + "." + $zhHash + ".zh-cn.xlf", "", "", $uuidA + "." + $zhHash + ".zh-cn.xlf")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("E3"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/" + $zhHandoffRepoCommitA + "/e2e/" + $uuidA + ".md", "", "", $uuidA + ".md")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("F3"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/" + $zhHandbackCommit + "/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/" + $uuidA + "." + $zhHash + ".zh-cn.xlf", "", "", $uuidA + "." + $zhHash + ".zh-cn.xlf")

$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/38beca1a811b8e315158bc4dda8d495b6acea141/.localization-config", "", "", ".localization-config")

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$deHandoffCommit = "a302c15b9158addb0facb1f663352e2d7ab956a9"
$deHandoffRepoCommitA = "7efcd252b5c61ef2b2755a4934c475b90a870a43"
$deHandbackCommit = "d011c6fc45bf6a5c2e398e946a5ddfd4fca9fe1c"
$deHash = "7d99b0cfb94d07c6c9898a5664a9f07e6b3c9ea5"
$deHashB = "fa8dfbcb6f034a590a88c3ecb061e2853a9c33f6"

# Row 2 now holds the 929aaf6c file (status/timestamps unchanged from before)
$wsDeDe.Range("A2").Value2 = ($uuidB + ".md")
$wsDeDe.Range("B2").Value2 = $handedBack
$wsDeDe.Range("C2").Value2 = ($uuidB + "." + $deHashB + ".de-de.xlf")
$wsDeDe.Range("D2").Value2 = "2016-03-08 18:49:32"
$wsDeDe.Range("E2").Value2 = ($uuidB + ".md")
$wsDeDe.Range("F2").Value2 = ($uuidB + "." + $deHashB + ".de-de.xlf")
$wsDeDe.Range("G2").Value2 = "2016-03-08 18:50:33"
$wsDeDe.Range("H2").Value2 = "Include"

# Row 3 now holds the 2acdd304 file, with refreshed status + handoff datetime
$wsDeDe.Range("A3").Value2 = ($uuidA + ".md")
$wsDeDe.Range("B3").Value2 = $readyForHandoff
$wsDeDe.Range("C3").Value2 = ($uuidA + "." + $deHash + ".de-de.xlf")
$wsDeDe.Range("D3").Value2 = "2016-03-08 18:51:18"
$wsDeDe.Range("E3").Value2 = ($uuidA + ".md")
$wsDeDe.Range("F3").Value2 = ($uuidA + "." + $deHash + ".de-de.xlf")
$wsDeDe.Range("G3").Value2 = "2016-03-08 18:50:33"
$wsDeDe.Range("H3").Value2 = "Include"

$wsDeDe.Hyperlinks.Delete()
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/38beca1a811b8e315158bc4dda8d495b6acea141/e2e/" + $uuidB + ".md", "", "", $uuidB + ".md")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/" + $deHandoffCommit + "/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/" + $uuidB + "." + $deHashB + ".de-de.xlf", "", "", $uuidB + "." + $deHashB + ".de-de.xlf")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("E2"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/" + $deHandoffRepoCommitA + "/e2e/" + $uuidB + ".md", "", "", $uuidB + ".md")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/" + $deHandbackCommit + "/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/" + $uuidB + "." + $deHashB + ".de-de.xlf", "", "", $uuidB + "." + $deHashB + ".de-de.xlf")

$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/38beca1a811b8e315158bc4dda8d495b6acea141/e2e/" + $uuidA + ".md", "", "", $uuidA + ".md")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("C3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/" + $deHandoffCommit + "/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/" + $uuidA + "." + $deHash + ".de-de.xlf", "", "", $uuidA + "." + $deHash + ".de-de.xlf")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("E3"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/" + $deHandoffRepoCommitA + "/e2e/" + $uuidA + ".md", "", "", $uuidA + ".md")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("F3"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/" + $deHandbackCommit + "/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/" + $uuidA + "." + $deHash + ".de-de.xlf", "", "", $uuidA + "." + $deHash + ".de-de.xlf")

$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/38beca1a811b8e315158bc4dda8d495b6acea141/.localization-config", "", "", ".localization-config")
